# Migration to Automation Org
# Update the QuantityBefore / QuantityAfter reference values on the
# "POReceipt" sheet to the new Automation-Org record numbers.
#
# These columns hold numeric-looking values that are stored as TEXT
# (e.g. "958.0"), so a plain Range.Value assignment would be re-interpreted
# by Excel as a number. Forcing the cell to Text format before the write
# (then clearing the format back, since the source cells carry no explicit
# style) preserves the original text representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POReceipt")

function Set-TextCellValue {
    param($Range, [string]$Value)

    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

# row -> (QuantityBefore "J", QuantityAfter "K")
$updates = @(
    @{ Row = 2; J = "1291.0"; K = "1295.0" },
    @{ Row = 3; J = "65.0";   K = "67.0"   },
    @{ Row = 4; J = "733.0";  K = "739.0"  }
)

foreach ($u in $updates) {
    Set-TextCellValue $ws.Range("J$($u.Row)") $u.J
    Set-TextCellValue $ws.Range("K$($u.Row)") $u.K
}
